$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04678466666666666
$ws.Range("H2").Value = 0.140354
$ws.Range("I2").Value = 0.006739448717762189
$ws.Range("J2").Value = 0.006739448717762188
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 3.606607652642222
$ws.Range("R2").Value = 32.45946887378
$ws.Range("S2").Value = 0.001620039919948361
$ws.Range("T2").Value = 0.001620039919948361
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04678466666666666
$ws.Range("H3").Value = 0.140354
$ws.Range("I3").Value = 0.006739448717762189
$ws.Range("J3").Value = 0.006739448717762188
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 4.752388186627555
$ws.Range("R3").Value = 42.771493679648
$ws.Range("S3").Value = 0.002134709211241003
$ws.Range("T3").Value = 0.002134709211241002
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04678466666666666
$ws.Range("H4").Value = 0.140354
$ws.Range("I4").Value = 0.006739448717762189
$ws.Range("J4").Value = 0.006739448717762188
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 6.644675996696889
$ws.Range("R4").Value = 59.80208397027199
$ws.Range("S4").Value = 0.002984699586572825
$ws.Range("T4").Value = 0.002984699586572824
$ws.Range("I5").Value = 0.9176013393810419
$ws.Range("J5").Value = 0.9176013393810418
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 491.0532227902022
$ws.Range("R5").Value = 4419.47900511182
$ws.Range("S5").Value = 0.2205745399438215
$ws.Range("T5").Value = 0.2205745399438215
$ws.Range("I6").Value = 0.9176013393810419
$ws.Range("J6").Value = 0.9176013393810418
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2906487033963527
$ws.Range("T6").Value = 0.2906487033963526
$ws.Range("I7").Value = 0.9176013393810419
$ws.Range("J7").Value = 0.9176013393810418
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.4063780960408678
$ws.Range("T7").Value = 0.4063780960408677
$ws.Range("G8").Value = 0.5252196666666668
$ws.Range("I8").Value = 0.07565921190119594
$ws.Range("J8").Value = 0.07565921190119593
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 40.48893374862556
$ws.Range("R8").Value = 364.40040373763
$ws.Range("S8").Value = 0.01818708750891258
$ws.Range("T8").Value = 0.01818708750891257
$ws.Range("G9").Value = 0.5252196666666668
$ws.Range("I9").Value = 0.07565921190119594
$ws.Range("J9").Value = 0.07565921190119593
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("Q9").Value = 53.3518333482009
$ws.Range("R9").Value = 480.1665001338081
$ws.Range("S9").Value = 0.02396492997046603
$ws.Range("T9").Value = 0.02396492997046602
$ws.Range("G10").Value = 0.5252196666666668
$ws.Range("I10").Value = 0.07565921190119594
$ws.Range("J10").Value = 0.07565921190119593
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("Q10").Value = 74.59526295139024
$ws.Range("R10").Value = 671.3573665625121
$ws.Range("S10").Value = 0.03350719442181734
$ws.Range("T10").Value = 0.03350719442181734
